$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (pushes existing data rows 2..113 down to 3..114)
$ws.Rows.Item(2).Insert()

# The inserted row copies formatting from the row above (the bold header row).
# Clear that so the new row matches the plain formatting of the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Re-apply the date number format used by the other rows' Fecha column (D).
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Fill in the new row's values (same Mercado/Región/etc. as every other row in this sheet).
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44643
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino ensalada"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 16000
$ws.Range("L2").Value = 17000
$ws.Range("M2").Value = 16571
$ws.Range("N2").Value = "$/caja 70 unidades"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 237
$ws.Range("Q2").Value = 70
$ws.Range("R2").Value = "Hortaliza"
